# Posner_cueing.xlsx - update targetX (column C) values and move the
# active selection, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C ("targetX") rows 2-11: the stimulus eccentricity was widened
# from +/-7 to +/-300.
$ws.Range("C2").Value = 300
$ws.Range("C3").Value = 300
$ws.Range("C4").Value = 300
$ws.Range("C5").Value = 300
$ws.Range("C6").Value = -300
$ws.Range("C7").Value = -300
$ws.Range("C8").Value = -300
$ws.Range("C9").Value = -300
$ws.Range("C10").Value = -300
$ws.Range("C11").Value = -300

# The author's last selection before saving moved to C12.
$ws.Range("C12").Select()
